# Requirements traceability sheet: add Sprint / Priority / QA evidence columns (C:L)
# Columns are filled one at a time (header, then its data) to mirror how the
# sheet was actually authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- C: Sprint number ----
$ws.Cells.Item(1, 3).Value = "Sprint"
for ($r = 2; $r -le 13; $r++)  { $ws.Cells.Item($r, 3).Value = 1 }
for ($r = 14; $r -le 25; $r++) { $ws.Cells.Item($r, 3).Value = 2 }
for ($r = 26; $r -le 37; $r++) { $ws.Cells.Item($r, 3).Value = 3 }
for ($r = 38; $r -le 49; $r++) { $ws.Cells.Item($r, 3).Value = 4 }

# ---- D: MoSCoW priority (only the first sprint's rows are filled in) ----
$ws.Cells.Item(1, 4).Value = "Priority (MoSCoW)"
for ($r = 2; $r -le 5; $r++) { $ws.Cells.Item($r, 4).Value = "Must" }

# ---- E:K single-cell headers for the remaining evidence columns ----
$ws.Cells.Item(1, 5).Value  = "Acceptance Criteria Summary"
$ws.Cells.Item(1, 6).Value  = "Code Modules / Files"
$ws.Cells.Item(1, 7).Value  = "Test Cases / Files"
$ws.Cells.Item(1, 8).Value  = "Coverage %"
$ws.Cells.Item(1, 9).Value  = "Status"
$ws.Cells.Item(1, 10).Value = "Evidence Link (Branch)"
$ws.Cells.Item(1, 11).Value = "Last Verified"

# L1 is left blank, but still needs to pick up the bold header style below.

# ---- Bold the whole new header row, C1:L1 ----
$ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item(1, 12)).Font.Bold = $true

# ---- Column widths (D:L) ----
$widths = @{
    4  = 18.5
    5  = 27.6666666666667
    6  = 22.5
    7  = 15.6666666666667
    8  = 11.6666666666667
    9  = 11.8333333333333
    10 = 20.6666666666667
    11 = 12.3333333333333
    12 = 13.3333333333333
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}

# ---- Selection, matching the saved workbook state ----
$ws.Range("E2").Select() | Out-Null
